# moved jcoin namespace fields back under custom property
#
# Column reshuffle on the "table-schema-baseline" sheet:
#   old col -> new col
#   A (jcoin:core_measure_section)  -> G (custom.jcoin:core_measure_section)
#   B (name)                         -> A (name)
#   C (type)                         -> B (type)
#   D (description)                  -> C (description)
#   E (trueValues)                   -> D (trueValues)
#   F (falseValues)                  -> E (falseValues)
#   G (constraints.enum)             -> F (constraints.enum)
#   H (title)                        -> I (title)
#   I (constraints.required)         -> J (constraints.required)
#   J (constraints.maxLength)        -> K (constraints.maxLength)
#   K (constraints.pattern)          -> L (constraints.pattern)
#   L (jcoin:final_variable_name)    -> H (custom.jcoin:final_variable_name)
#   M (format)                       -> M (format)
#   N (jcoin:notes)                  -> N (custom.jcoin:notes)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table-schema-baseline")

$firstRow = 1
$lastRow = 38
$lastCol = 14

$srcRange = $ws.Range("A1:N38")
$src = $srcRange.Value2
$dst = $ws.Range("A1:N38").Value2

# old column index (1-based) -> new column index (1-based)
$colMap = @{ 1 = 7; 2 = 1; 3 = 2; 4 = 3; 5 = 4; 6 = 5; 7 = 6; 8 = 9; 9 = 10; 10 = 11; 11 = 12; 12 = 8; 13 = 13; 14 = 14 }

for ($r = 1; $r -le $lastRow; $r++) {
    for ($oldC = 1; $oldC -le $lastCol; $oldC++) {
        $newC = $colMap[$oldC]
        $dst[$r, $newC] = $src[$r, $oldC]
    }
}

# Rename the three header cells that gained the "custom." prefix
$dst[1, 7] = "custom.jcoin:core_measure_section"
$dst[1, 8] = "custom.jcoin:final_variable_name"
$dst[1, 14] = "custom.jcoin:notes"

$ws.Range("A1:N38").Value2 = $dst

# New column widths (stored/display width). ColumnWidth property is offset
# from the stored width by 5/6 (5px padding / 6px max-digit-width), so
# subtract that to land exactly on the target stored width.
$offset = 5 / 6
$widths = @(24, 10, 76, 10, 11, 164, 33, 32, 36.5, 20, 21, 41, 6, 171)
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c - 1] - $offset
}
